$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1 (18:22 -> 18:52)
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 18:52"

# Row 4
$ws.Cells.Item(4, 2).Value = 478366
$ws.Cells.Item(4, 3).Value = 9800
$ws.Cells.Item(4, 4).Value = 26163
$ws.Cells.Item(4, 5).Value = 434276
$ws.Cells.Item(4, 6).Value = 10896
$ws.Cells.Item(4, 7).Value = 1236
$ws.Cells.Item(4, 8).Value = 17927

# Row 10
$ws.Cells.Item(10, 5).Value = 61179
$ws.Cells.Item(10, 7).Value = 980
$ws.Cells.Item(10, 8).Value = 8958

# Row 12
$ws.Cells.Item(12, 2).Value = 47029
$ws.Cells.Item(12, 3).Value = 4747
$ws.Cells.Item(12, 4).Value = 2423
$ws.Cells.Item(12, 5).Value = 43600
$ws.Cells.Item(12, 6).Value = 1667
$ws.Cells.Item(12, 7).Value = 98
$ws.Cells.Item(12, 8).Value = 1006

# Row 23
$ws.Cells.Item(23, 4).Value = 381
$ws.Cells.Item(23, 5).Value = 8434

# Row 25
$ws.Cells.Item(25, 1).Value = "Ecuador"
$ws.Cells.Item(25, 2).Value = 7161
$ws.Cells.Item(25, 3).Value = 2196
$ws.Cells.Item(25, 4).Value = 368
$ws.Cells.Item(25, 5).Value = 6496
$ws.Cells.Item(25, 6).Value = 171
$ws.Cells.Item(25, 7).Value = 25
$ws.Cells.Item(25, 8).Value = 297

# Row 26
$ws.Cells.Item(26, 1).Value = "Irlanda"
$ws.Cells.Item(26, 2).Value = 6574
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 25
$ws.Cells.Item(26, 5).Value = 6286
$ws.Cells.Item(26, 6).Value = 194
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 263

# Row 27
$ws.Cells.Item(27, 1).Value = "Chile"
$ws.Cells.Item(27, 2).Value = 6501
$ws.Cells.Item(27, 3).Value = 529
$ws.Cells.Item(27, 4).Value = 1571
$ws.Cells.Item(27, 5).Value = 4865
$ws.Cells.Item(27, 6).Value = 360
$ws.Cells.Item(27, 7).Value = 8
$ws.Cells.Item(27, 8).Value = 65

# Row 28
$ws.Cells.Item(28, 1).Value = "Noruega"
$ws.Cells.Item(28, 2).Value = 6244
$ws.Cells.Item(28, 3).Value = 25
$ws.Cells.Item(28, 4).Value = 32
$ws.Cells.Item(28, 5).Value = 6100
$ws.Cells.Item(28, 6).Value = 70
$ws.Cells.Item(28, 7).Value = 4
$ws.Cells.Item(28, 8).Value = 112

# Row 29
$ws.Cells.Item(29, 1).Value = "Australia"
$ws.Cells.Item(29, 2).Value = 6203
$ws.Cells.Item(29, 3).Value = 51
$ws.Cells.Item(29, 4).Value = 3141
$ws.Cells.Item(29, 5).Value = 3009
$ws.Cells.Item(29, 6).Value = 74
$ws.Cells.Item(29, 7).Value = 2
$ws.Cells.Item(29, 8).Value = 53

# Row 30
$ws.Cells.Item(30, 1).Value = "Polonia"
$ws.Cells.Item(30, 2).Value = 5955
$ws.Cells.Item(30, 3).Value = 380
$ws.Cells.Item(30, 4).Value = 318
$ws.Cells.Item(30, 5).Value = 5456
$ws.Cells.Item(30, 6).Value = 160
$ws.Cells.Item(30, 7).Value = 7
$ws.Cells.Item(30, 8).Value = 181

# Row 31
$ws.Cells.Item(31, 1).Value = "Dinamarca"
$ws.Cells.Item(31, 2).Value = 5819
$ws.Cells.Item(31, 3).Value = 184
$ws.Cells.Item(31, 4).Value = 1773
$ws.Cells.Item(31, 5).Value = 3799
$ws.Cells.Item(31, 6).Value = 113
$ws.Cells.Item(31, 7).Value = 10
$ws.Cells.Item(31, 8).Value = 247

# Row 32
$ws.Cells.Item(32, 1).Value = "Chequia"
$ws.Cells.Item(32, 2).Value = 5674
$ws.Cells.Item(32, 3).Value = 105
$ws.Cells.Item(32, 4).Value = 346
$ws.Cells.Item(32, 5).Value = 5209
$ws.Cells.Item(32, 6).Value = 98
$ws.Cells.Item(32, 7).Value = 7
$ws.Cells.Item(32, 8).Value = 119

# Row 33
$ws.Cells.Item(33, 1).Value = "Japon"
$ws.Cells.Item(33, 2).Value = 5530
$ws.Cells.Item(33, 3).Value = 183
$ws.Cells.Item(33, 4).Value = 685
$ws.Cells.Item(33, 5).Value = 4746
$ws.Cells.Item(33, 6).Value = 109
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 99

# Row 34
$ws.Cells.Item(34, 1).Value = "Rumania"
$ws.Cells.Item(34, 2).Value = 5467
$ws.Cells.Item(34, 3).Value = 265
$ws.Cells.Item(34, 4).Value = 729
$ws.Cells.Item(34, 5).Value = 4468
$ws.Cells.Item(34, 6).Value = 183
$ws.Cells.Item(34, 7).Value = 22
$ws.Cells.Item(34, 8).Value = 270

# Row 35
$ws.Cells.Item(35, 1).Value = "Peru"
$ws.Cells.Item(35, 2).Value = 5256
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 1438
$ws.Cells.Item(35, 5).Value = 3680
$ws.Cells.Item(35, 6).Value = 124
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 138

# Row 75
$ws.Cells.Item(75, 2).Value = 812
$ws.Cells.Item(75, 3).Value = 31
$ws.Cells.Item(75, 5).Value = 738

# Row 86
$ws.Cells.Item(86, 6).Value = 11
